# "review of the paper, major revisions" — refresh the cached survey
# statistics on the two raw-data sheets ("exit" and "basal"); the
# balance_response summary sheet recalculates its ROUND()/CONCATENATE()
# formulas automatically off these inputs.

$wb = $excel.ActiveWorkbook
$wsExit = $wb.Worksheets.Item("exit")
$wsBasal = $wb.Worksheets.Item("basal")

# --- exit sheet updates ---
$wsExit.Range("B2").Value = 2157.521888206682
$wsExit.Range("D2").Value = 2171.6630968317713
$wsExit.Range("E2").Value = 0.17251166970702933
$wsExit.Range("B3").Value = 30.911295389306172
$wsExit.Range("D3").Value = 31.80489829244328
$wsExit.Range("B4").Value = 0.18236185312176062
$wsExit.Range("D4").Value = 0.18310278149635578
$wsExit.Range("E4").Value = 0.71165592726145865
$wsExit.Range("B5").Value = 0.024502662183247662
$wsExit.Range("D5").Value = 0.024355796665085788
$wsExit.Range("B6").Value = 32.565573770491802
$wsExit.Range("C6").Value = 37.243902439024389
$wsExit.Range("D6").Value = 33.036855036855037
$wsExit.Range("E6").Value = 0.28366897076976583
$wsExit.Range("B7").Value = 1.0480981897354886
$wsExit.Range("C7").Value = 4.2716850630413781
$wsExit.Range("D7").Value = 1.036376323398289
$wsExit.Range("B8").Value = 12541
$wsExit.Range("D8").Value = 13446
$wsExit.Range("E3").ClearContents()
$wsExit.Range("E5").ClearContents()
$wsExit.Range("E7").ClearContents()
$wsExit.Range("E8").ClearContents()

# --- basal sheet updates ---
$wsBasal.Range("C2").Value = 2145.5996742983043
$wsBasal.Range("D2").Value = 2171.6630968317713
$wsBasal.Range("E2").Value = 0.077673504954345954
$wsBasal.Range("C3").Value = 35.438239354568644
$wsBasal.Range("D3").Value = 31.804898292443283
$wsBasal.Range("C4").Value = 0.18344669029600535
$wsBasal.Range("D4").Value = 0.18310278149635578
$wsBasal.Range("E4").Value = 0.94749239594595602
$wsBasal.Range("C5").Value = 0.024514901266524693
$wsBasal.Range("D5").Value = 0.024355796665085788
$wsBasal.Range("B6").Value = 34.938271604938272
$wsBasal.Range("C6").Value = 32.564417177914109
$wsBasal.Range("D6").Value = 33.036855036855037
$wsBasal.Range("E6").Value = 0.41501689580331746
$wsBasal.Range("B7").Value = 2.7001569272532908
$wsBasal.Range("C7").Value = 1.1071527165208588
$wsBasal.Range("D7").Value = 1.036376323398289
$wsBasal.Range("C8").Value = 10439
$wsBasal.Range("D8").Value = 13446
$wsBasal.Range("E3").ClearContents()
$wsBasal.Range("E5").ClearContents()
$wsBasal.Range("E7").ClearContents()
$wsBasal.Range("E8").ClearContents()
